$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-17.03079643509903, 0.134081037299484, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903),
    @(-17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, 1.188603728550447, -17.03079643509903, -17.03079643509903),
    @(-17.03079643509903, 0.1918192515794448, 1.062419938535544, -17.03079643509903, 3.864557698741737, -17.03079643509903, 1.659300461680259, -17.03079643509903, 2.439729388308265, -17.03079643509903),
    @(-17.03079643509903, 0.5418815895332343, -17.03079643509903, -17.03079643509903, -17.03079643509903, 3.516227040956146, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903),
    @(-17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903),
    @(2.567294240711527, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903),
    @(-17.03079643509903, -17.03079643509903, -17.03079643509903, 4.321917859150612, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903),
    @(3.814847512667724, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903),
    @(-17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, 0.712019550245869, -17.03079643509903, 1.726805175899547),
    @(-17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, 1.59869181692837, -17.03079643509903, -17.03079643509903, -17.03079643509903, 1.205521463876104),
    @(-17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903),
    @(-17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, 1.136234522543681, 2.272344378294263),
    @(-17.03079643509903, -17.03079643509903, 3.299014134181536, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, 1.816043672893322),
    @(-17.03079643509903, -17.03079643509903, 0.07454943380528492, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903),
    @(-17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, 2.625651123742962, -17.03079643509903),
    @(-17.03079643509903, 0.5929207652909531, 0.3021928978728465, -17.03079643509903, -17.03079643509903, -17.03079643509903, 2.706379925607641, 0.4045564110337567, 1.777592763757902, -17.03079643509903),
    @(-17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, -17.03079643509903, 1.911010592280806, 0.04247695804470942, 1.473141853807686, -17.03079643509903),
    @(-17.03079643509903, -17.03079643509903, 1.702314905893077, -17.03079643509903, -17.03079643509903, -17.03079643509903, 1.427486006942833, 2.205519127911828, -17.03079643509903, -17.03079643509903),
    @(-17.03079643509903, 1.665223604065123, 1.338430273330074, -17.03079643509903, 2.441912250531439, -17.03079643509903, 0.9960133480898876, 3.188482002328875, -17.03079643509903, 2.592522463877847),
    @(-17.03079643509903, 3.539097175466206, -17.03079643509903, -17.03079643509903, -17.03079643509903, 2.467150883079579, 0.9026233664207445, -17.03079643509903, -17.03079643509903, -17.03079643509903)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}

Write-Output "done"